{"js": "// Office.js (Word JavaScript API) implementation of the target edit:\n//\n//   {x3}{g4}ommentscay{/g5}{x6} ... {x7}{g8}.{/g9}\n// becomes\n//   {x3}ommentscay{x4} ... {x5}.\n//\n// i.e. the opening/closing placeholder markers ({g4}/{/g5} and\n// {g8}/{/g9}) are removed (merging the surrounding runs), and the\n// numbered placeholders that followed the removed comment text are\n// renumbered down by two ({x6} -> {x4}, {x7} -> {x5}).\n\nasync function findUnique(text) {\n  const results = context.document.body.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  return results.items[0];\n}\n\n// 1) Drop the opening marker right after \"{x3}\" so the two runs\n//    \"{x3}\" and \"{g4}ommentscay{/g5}\" collapse into \"{x3}ommentscay\".\n(await findUnique(\"{g4}\")).insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 2) Drop the matching closing marker.\n(await findUnique(\"{/g5}\")).insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 3) Renumber the placeholder that sits right after the comment range end.\n(await findUnique(\"{x6}\")).insertText(\"{x4}\", \"Replace\");\nawait context.sync();\n\n// 4) Drop the second opening marker, collapsing \"{x7}\" and\n//    \"{g8}.{/g9}\" into \"{x7}.\".\n(await findUnique(\"{g8}\")).insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 5) Drop the matching closing marker.\n(await findUnique(\"{/g9}\")).insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 6) Renumber the final placeholder.\n(await findUnique(\"{x7}\")).insertText(\"{x5}\", \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) implementation of the target edit:\n#\n#   {x3}{g4}ommentscay{/g5}{x6} ... {x7}{g8}.{/g9}\n# becomes\n#   {x3}ommentscay{x4} ... {x5}.\n#\n# i.e. the opening/closing placeholder markers ({g4}/{/g5} and\n# {g8}/{/g9}) are removed (merging the surrounding runs), and the\n# numbered placeholders that followed the removed comment text are\n# renumbered down by two ({x6} -> {x4}, {x7} -> {x5}).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Unique([string]$needle, [string]$replacement) {\n    # $needle is unique in the document, so a single Find.Execute() locates\n    # exactly the run we want to touch; plain Range.Text assignment is fine\n    # here because neither side of the match sits right against a\n    # zero-width marker (commentRangeEnd / commentReference).\n    $rng = $d.Content\n    $rng.Find.Text = $needle\n    $rng.Find.Execute() | Out-Null\n    $rng.Text = $replacement\n}\n\nfunction Rename-AfterMarker([string]$needle, [string]$replacement) {\n    # $needle immediately follows a zero-width marker (commentRangeEnd or\n    # the comment-reference run). Assigning straight to Range.Text in that\n    # spot relocates the marker to the wrong side of the new text, so\n    # instead insert the replacement right after the match, then delete\n    # the original text in a separate pass.\n    $rng = $d.Content\n    $rng.Find.Text = $needle\n    $rng.Find.Execute() | Out-Null\n    $rng.InsertAfter($replacement)\n\n    $rng2 = $d.Content\n    $rng2.Find.Text = $needle\n    $rng2.Find.Execute() | Out-Null\n    $rng2.Delete()\n}\n\n# 1) Remove the opening/closing markers wrapping \"ommentscay\" so \"{x3}\" and\n#    \"{g4}ommentscay{/g5}\" collapse into a single \"{x3}ommentscay\" run.\nReplace-Unique \"{g4}\" \"\"\nReplace-Unique \"{/g5}\" \"\"\n\n# 2) Renumber the placeholder right after the comment range end.\nRename-AfterMarker \"{x6}\" \"{x4}\"\n\n# 3) Remove the opening/closing markers wrapping \".\" so \"{x7}\" and\n#    \"{g8}.{/g9}\" collapse into a single \"{x7}.\" run.\nReplace-Unique \"{g8}\" \"\"\nReplace-Unique \"{/g9}\" \"\"\n\n# 4) Renumber the final placeholder (right after the comment reference run).\nRename-AfterMarker \"{x7}\" \"{x5}\"\n"}
